$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("buffer")

# The "buffer" sheet used to have a "name" column (A) holding computed
# labels like "product @ factory" / "product @ store". That column is
# replaced by using the buffer's id (location + item) instead, so we
# simply delete the old "name" column, shifting location/item left.
$ws.Columns.Item(1).Delete()
